$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the numeric-looking strings (prices & most holiday dates)
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("B2:B6").NumberFormat = "@"

# Ticket prices (column A, rows 2-5)
$ws.Range("A2").Value = "10"
$ws.Range("A3").Value = "20"
$ws.Range("A4").Value = "30"
$ws.Range("A5").Value = "40"

# Holiday dates (column B, rows 2-6 stored as text, row 7 stored as a number)
$ws.Range("B2").Value = "1012022"
$ws.Range("B3").Value = "1022022"
$ws.Range("B4").Value = "2012022"
$ws.Range("B5").Value = "2022022"
$ws.Range("B6").Value = "9082022"
$ws.Range("B7").Value = 6012022

# Restore the default ("Normal") cell style so the text formatting we used
# only to coerce these into text values doesn't stick around as a new style.
$ws.Range("A2:A5").Style = "Normal"
$ws.Range("B2:B7").Style = "Normal"

$ws.Range("B8").Select() | Out-Null
